$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = (Get-Date -Year 2021 -Month 7 -Day 13).Date
$ws.Cells.Item(2, 10).Value = 50

# Row 3
$ws.Cells.Item(3, 4).Value = (Get-Date -Year 2021 -Month 4 -Day 6).Date
$ws.Cells.Item(3, 10).Value = 40

# Row 4
$ws.Cells.Item(4, 4).Value = (Get-Date -Year 2020 -Month 12 -Day 1).Date
$ws.Cells.Item(4, 10).Value = 45
$ws.Cells.Item(4, 12).Value = 2500
$ws.Cells.Item(4, 13).Value = 2500
$ws.Cells.Item(4, 16).Value = 833

# Row 5
$ws.Cells.Item(5, 4).Value = (Get-Date -Year 2023 -Month 3 -Day 7).Date
$ws.Cells.Item(5, 10).Value = 45
$ws.Cells.Item(5, 11).Value = 4000
$ws.Cells.Item(5, 12).Value = 4000
$ws.Cells.Item(5, 13).Value = 4000
$ws.Cells.Item(5, 16).Value = 1333

# Row 6
$ws.Cells.Item(6, 4).Value = (Get-Date -Year 2020 -Month 11 -Day 30).Date
$ws.Cells.Item(6, 10).Value = 68
$ws.Cells.Item(6, 12).Value = 3000
$ws.Cells.Item(6, 13).Value = 3000
$ws.Cells.Item(6, 16).Value = 1000

# Row 7
$ws.Cells.Item(7, 4).Value = (Get-Date -Year 2020 -Month 12 -Day 14).Date
$ws.Cells.Item(7, 10).Value = 78
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 13).Value = 3000
$ws.Cells.Item(7, 16).Value = 1000

# Row 8
$ws.Cells.Item(8, 4).Value = (Get-Date -Year 2021 -Month 1 -Day 26).Date

# Row 9
$ws.Cells.Item(9, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 29).Date
$ws.Cells.Item(9, 10).Value = 68
$ws.Cells.Item(9, 11).Value = 2000
$ws.Cells.Item(9, 12).Value = 2000
$ws.Cells.Item(9, 13).Value = 2000
$ws.Cells.Item(9, 16).Value = 667

# Row 10
$ws.Cells.Item(10, 4).Value = (Get-Date -Year 2020 -Month 12 -Day 22).Date
$ws.Cells.Item(10, 10).Value = 65
$ws.Cells.Item(10, 11).Value = 3000
$ws.Cells.Item(10, 13).Value = 3000
$ws.Cells.Item(10, 16).Value = 1000

# Row 11
$ws.Cells.Item(11, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 14).Date
$ws.Cells.Item(11, 10).Value = 104
$ws.Cells.Item(11, 11).Value = 2800
$ws.Cells.Item(11, 13).Value = 2904
$ws.Cells.Item(11, 16).Value = 968

# Row 12
$ws.Cells.Item(12, 4).Value = (Get-Date -Year 2023 -Month 6 -Day 15).Date
$ws.Cells.Item(12, 10).Value = 90
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 12).Value = 3500
$ws.Cells.Item(12, 13).Value = 3278
$ws.Cells.Item(12, 16).Value = 1093

# Row 13
$ws.Cells.Item(13, 4).Value = (Get-Date -Year 2021 -Month 2 -Day 16).Date
$ws.Cells.Item(13, 10).Value = 45
$ws.Cells.Item(13, 11).Value = 3000
$ws.Cells.Item(13, 13).Value = 3000
$ws.Cells.Item(13, 16).Value = 1000

# Row 14
$ws.Cells.Item(14, 4).Value = (Get-Date -Year 2021 -Month 4 -Day 5).Date
$ws.Cells.Item(14, 10).Value = 45
$ws.Cells.Item(14, 12).Value = 3000
$ws.Cells.Item(14, 13).Value = 3000
$ws.Cells.Item(14, 16).Value = 1000

# Row 15
$ws.Cells.Item(15, 4).Value = (Get-Date -Year 2023 -Month 2 -Day 8).Date
$ws.Cells.Item(15, 10).Value = 87

# Row 16
$ws.Cells.Item(16, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 13).Date
$ws.Cells.Item(16, 10).Value = 50

# Row 17
$ws.Cells.Item(17, 4).Value = (Get-Date -Year 2021 -Month 1 -Day 29).Date
$ws.Cells.Item(17, 10).Value = 56

# Row 18
$ws.Cells.Item(18, 4).Value = (Get-Date -Year 2023 -Month 1 -Day 9).Date
$ws.Cells.Item(18, 10).Value = 78
$ws.Cells.Item(18, 11).Value = 3000
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 13).Value = 3000
$ws.Cells.Item(18, 16).Value = 1000

# Row 19
$ws.Cells.Item(19, 4).Value = (Get-Date -Year 2021 -Month 2 -Day 15).Date
$ws.Cells.Item(19, 10).Value = 95
$ws.Cells.Item(19, 11).Value = 2500
$ws.Cells.Item(19, 13).Value = 2737
$ws.Cells.Item(19, 16).Value = 912

# Row 20
$ws.Cells.Item(20, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 27).Date
$ws.Cells.Item(20, 10).Value = 104
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = 2260
$ws.Cells.Item(20, 16).Value = 753

# Row 21
$ws.Cells.Item(21, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 6).Date
$ws.Cells.Item(21, 10).Value = 125
$ws.Cells.Item(21, 11).Value = 2200
$ws.Cells.Item(21, 12).Value = 2200
$ws.Cells.Item(21, 13).Value = 2200
$ws.Cells.Item(21, 16).Value = 733

# Row 22
$ws.Cells.Item(22, 4).Value = (Get-Date -Year 2023 -Month 1 -Day 11).Date
$ws.Cells.Item(22, 10).Value = 68
$ws.Cells.Item(22, 11).Value = 3500
$ws.Cells.Item(22, 12).Value = 3500
$ws.Cells.Item(22, 13).Value = 3500
$ws.Cells.Item(22, 16).Value = 1167

# Row 23
$ws.Cells.Item(23, 4).Value = (Get-Date -Year 2021 -Month 5 -Day 24).Date
$ws.Cells.Item(23, 10).Value = 54

# Row 24
$ws.Cells.Item(24, 4).Value = (Get-Date -Year 2020 -Month 12 -Day 28).Date
$ws.Cells.Item(24, 10).Value = 70
$ws.Cells.Item(24, 11).Value = 3000
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 13).Value = 3000
$ws.Cells.Item(24, 16).Value = 1000

# Row 25
$ws.Cells.Item(25, 4).Value = (Get-Date -Year 2022 -Month 4 -Day 18).Date
$ws.Cells.Item(25, 10).Value = 92
$ws.Cells.Item(25, 11).Value = 2500
$ws.Cells.Item(25, 13).Value = 2755
$ws.Cells.Item(25, 16).Value = 918

# Row 26
$ws.Cells.Item(26, 4).Value = (Get-Date -Year 2021 -Month 7 -Day 12).Date
$ws.Cells.Item(26, 10).Value = 81
$ws.Cells.Item(26, 11).Value = 2800
$ws.Cells.Item(26, 12).Value = 3000
$ws.Cells.Item(26, 13).Value = 2889
$ws.Cells.Item(26, 16).Value = 963

# Row 27
$ws.Cells.Item(27, 4).Value = (Get-Date -Year 2022 -Month 10 -Day 11).Date
$ws.Cells.Item(27, 10).Value = 80

# Row 28
$ws.Cells.Item(28, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 7).Date
$ws.Cells.Item(28, 10).Value = 88
$ws.Cells.Item(28, 11).Value = 2000
$ws.Cells.Item(28, 12).Value = 2200
$ws.Cells.Item(28, 13).Value = 2091
$ws.Cells.Item(28, 16).Value = 697

# Row 29
$ws.Cells.Item(29, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 31).Date
$ws.Cells.Item(29, 10).Value = 85
$ws.Cells.Item(29, 11).Value = 3000
$ws.Cells.Item(29, 12).Value = 3000
$ws.Cells.Item(29, 13).Value = 3000
$ws.Cells.Item(29, 16).Value = 1000

# Row 30
$ws.Cells.Item(30, 4).Value = (Get-Date -Year 2021 -Month 1 -Day 28).Date
$ws.Cells.Item(30, 10).Value = 67
$ws.Cells.Item(30, 11).Value = 3000
$ws.Cells.Item(30, 12).Value = 3000
$ws.Cells.Item(30, 13).Value = 3000
$ws.Cells.Item(30, 16).Value = 1000

# Row 31
$ws.Cells.Item(31, 4).Value = (Get-Date -Year 2021 -Month 1 -Day 27).Date
$ws.Cells.Item(31, 10).Value = 80
$ws.Cells.Item(31, 12).Value = 3000
$ws.Cells.Item(31, 13).Value = 2781
$ws.Cells.Item(31, 16).Value = 927

# Row 32
$ws.Cells.Item(32, 4).Value = (Get-Date -Year 2021 -Month 3 -Day 5).Date
$ws.Cells.Item(32, 10).Value = 60
$ws.Cells.Item(32, 11).Value = 3500
$ws.Cells.Item(32, 12).Value = 3500
$ws.Cells.Item(32, 13).Value = 3500
$ws.Cells.Item(32, 16).Value = 1167

# Row 33
$ws.Cells.Item(33, 4).Value = (Get-Date -Year 2023 -Month 2 -Day 10).Date
$ws.Cells.Item(33, 10).Value = 110
$ws.Cells.Item(33, 12).Value = 3300
$ws.Cells.Item(33, 13).Value = 3136
$ws.Cells.Item(33, 16).Value = 1045

# Row 34
$ws.Cells.Item(34, 4).Value = (Get-Date -Year 2021 -Month 1 -Day 25).Date
$ws.Cells.Item(34, 10).Value = 50
$ws.Cells.Item(34, 11).Value = 2500
$ws.Cells.Item(34, 12).Value = 2500
$ws.Cells.Item(34, 13).Value = 2500
$ws.Cells.Item(34, 16).Value = 833

# Row 35
$ws.Cells.Item(35, 4).Value = (Get-Date -Year 2022 -Month 3 -Day 7).Date
$ws.Cells.Item(35, 10).Value = 78
$ws.Cells.Item(35, 11).Value = 3500
$ws.Cells.Item(35, 12).Value = 3500
$ws.Cells.Item(35, 13).Value = 3500
$ws.Cells.Item(35, 16).Value = 1167
